$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output ("WARNING: could not find: " + $find)
    }
    return $ok
}

# 1) "flex-box" -> "flex"
Replace-Text "flex-box" "flex"

# 2) fix "such ass ability" typo -> "such as ability"
Replace-Text "such ass ability" "such as ability"

# 3) "functionality will happen" -> "functionalities will happen"
Replace-Text "functionality will happen" "functionalities will happen"

# 4) Replace the "Ranks and bonuses" paragraph with three new paragraphs describing
#    the skill table's visual appearance.
$target = "The Ranks and bonuses are entered by the user"
$rng = $d.Content
$found = $rng.Find.Execute($target)
if ($found) {
    $para = $rng.Paragraphs(1)
    $pRng = $para.Range
    # trim the trailing paragraph mark so we only replace the paragraph's content
    $pRng.End = $pRng.End - 1
    $pRng.Text = ""

    $t1 = "For the visual appearance, the skill list and its many fields will be represented as a table. The headers will be tilted 45 degrees clockwise, pivoting from the top of the column. This is to minimize the overall width of the table. "
    $t2 = "The skill name column is generated from the skills table in the database, this allows for any custom skills to be added to the list by a player or dungeon master. The Untrained column will be marked with a special character from the UTF-8 table such as a filled square. Next is the Skill Bonus column which shows aa calculated total from the remaining columns using the above mathematical formula."
    $t3 = "The Class Skill column uses the same special character as the untrained column to denote that skill as a class skill for those characters chosen class. The Armor Check Penalty column (abbreviated ACP) provides a negative penalty to that skill based on the ACP of any armor worn. The Ranks, Racial, Feats, and Misc columns are numeric form fields entered by the user and MUST contain an integer value."

    $pRng.InsertAfter($t1 + "`r" + $t2 + "`r" + $t3)
} else {
    Write-Output "WARNING: could not find Ranks and bonuses paragraph"
}

Write-Output "done"
